# Update the "Flight ID" column (column A) on the OutAssignment sheet.
# Rows 23-43 change from 13 -> 5
# Rows 44-52 change from 5  -> 6
# Rows 53-57 change from 6  -> 7
# Rows 58-64 change from 7  -> 8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OutAssignment")

for ($r = 23; $r -le 43; $r++) {
    $ws.Cells.Item($r, 1).Value = 5
}

for ($r = 44; $r -le 52; $r++) {
    $ws.Cells.Item($r, 1).Value = 6
}

for ($r = 53; $r -le 57; $r++) {
    $ws.Cells.Item($r, 1).Value = 7
}

for ($r = 58; $r -le 64; $r++) {
    $ws.Cells.Item($r, 1).Value = 8
}
